$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1614
$ws.Cells.Item(29, 10).Value = 7950
$ws.Cells.Item(29, 12).Value = 23850
$ws.Cells.Item(29, 14).Value = -24412

$ws.Cells.Item(31, 8).Value = 1265101.1
$ws.Cells.Item(31, 9).Value = 1265101.1
$ws.Cells.Item(31, 11).Value = 3795303.3
$ws.Cells.Item(31, 13).Value = -3795073.3

$ws.Cells.Item(38, 8).Value = 1298.2727
$ws.Cells.Item(38, 9).Value = 86.77778
$ws.Cells.Item(38, 10).Value = 6750
$ws.Cells.Item(38, 11).Value = 260.33334
$ws.Cells.Item(38, 12).Value = 20250
$ws.Cells.Item(38, 13).Value = 111.66666
$ws.Cells.Item(38, 14).Value = -20994

$ws.Cells.Item(58, 8).Value = 985.35
$ws.Cells.Item(58, 9).Value = 171.3077
$ws.Cells.Item(58, 10).Value = 2497.1428
$ws.Cells.Item(58, 11).Value = 513.9231
$ws.Cells.Item(58, 12).Value = 7491.428400000001
$ws.Cells.Item(58, 13).Value = -363.9231
$ws.Cells.Item(58, 14).Value = -7791.428400000001

$ws.Cells.Item(70, 8).Value = 15243652
$ws.Cells.Item(70, 9).Value = 41917492
$ws.Cells.Item(70, 10).Value = 1457.1428
$ws.Cells.Item(70, 11).Value = 125752476
$ws.Cells.Item(70, 12).Value = 4371.428400000001
$ws.Cells.Item(70, 13).Value = -125752206
$ws.Cells.Item(70, 14).Value = -4911.428400000001

$ws.Cells.Item(73, 8).Value = 15243652
$ws.Cells.Item(73, 9).Value = 41917492
$ws.Cells.Item(73, 10).Value = 1457.1428
$ws.Cells.Item(73, 11).Value = 125752476
$ws.Cells.Item(73, 12).Value = 4371.428400000001
$ws.Cells.Item(73, 13).Value = -125751540
$ws.Cells.Item(73, 14).Value = -6243.428400000001

$ws.Cells.Item(88, 8).Value = 2000
$ws.Cells.Item(88, 10).Value = 500
$ws.Cells.Item(88, 12).Value = 500
$ws.Cells.Item(88, 14).Value = -1312

$ws.Cells.Item(91, 8).Value = 2000
$ws.Cells.Item(91, 10).Value = 500
$ws.Cells.Item(91, 12).Value = 500
$ws.Cells.Item(91, 14).Value = -3308

$ws.Cells.Item(92, 8).Value = 1099.8572
$ws.Cells.Item(92, 9).Value = 1195.4445
$ws.Cells.Item(92, 10).Value = 526.3333
$ws.Cells.Item(92, 11).Value = 1195.4445
$ws.Cells.Item(92, 12).Value = 526.3333
$ws.Cells.Item(92, 13).Value = 52.55549999999994
$ws.Cells.Item(92, 14).Value = -3022.3333

$ws.Cells.Item(107, 8).Value = 492.85715
$ws.Cells.Item(107, 9).Value = 475
$ws.Cells.Item(107, 10).Value = 516.6667
$ws.Cells.Item(107, 11).Value = 475
$ws.Cells.Item(107, 12).Value = 516.6667
$ws.Cells.Item(107, 13).Value = 1445
$ws.Cells.Item(107, 14).Value = -4356.6667

$ws.Cells.Item(116, 8).Value = 15387057
$ws.Cells.Item(116, 9).Value = 50001600
$ws.Cells.Item(116, 10).Value = 2815.7778
$ws.Cells.Item(116, 11).Value = 50001600
$ws.Cells.Item(116, 12).Value = 2815.7778
$ws.Cells.Item(116, 13).Value = -49998158
$ws.Cells.Item(116, 14).Value = -9699.7778

$ws.Cells.Item(131, 8).Value = 6189.8667
$ws.Cells.Item(131, 9).Value = 1268.1538
$ws.Cells.Item(131, 10).Value = 9953.529
$ws.Cells.Item(131, 11).Value = 3804.4614
$ws.Cells.Item(131, 12).Value = 29860.587
$ws.Cells.Item(131, 13).Value = 1235.5386
$ws.Cells.Item(131, 14).Value = -39940.587

$ws.Cells.Item(137, 8).Value = 2024.5186
$ws.Cells.Item(137, 9).Value = 1315.1904
$ws.Cells.Item(137, 10).Value = 4507.1665
$ws.Cells.Item(137, 11).Value = 3945.5712
$ws.Cells.Item(137, 12).Value = 13521.4995
$ws.Cells.Item(137, 13).Value = -1395.5712
$ws.Cells.Item(137, 14).Value = -18621.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 3000
$ws.Cells.Item(27, 9).Value = 3000
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 3000
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -2816

$ws.Cells.Item(61, 8).Value = 3039.8262
$ws.Cells.Item(61, 9).Value = 2651
$ws.Cells.Item(61, 10).Value = 3928.5715
$ws.Cells.Item(61, 11).Value = 2651
$ws.Cells.Item(61, 12).Value = 3928.5715
$ws.Cells.Item(61, 13).Value = -2439
$ws.Cells.Item(61, 14).Value = -4352.5715

$ws.Cells.Item(88, 8).Value = 3021.25
$ws.Cells.Item(88, 9).Value = 2400
$ws.Cells.Item(88, 10).Value = 3228.3333
$ws.Cells.Item(88, 11).Value = 2400
$ws.Cells.Item(88, 12).Value = 3228.3333
$ws.Cells.Item(88, 13).Value = -1994
$ws.Cells.Item(88, 14).Value = -4040.3333

$ws.Cells.Item(91, 8).Value = 3021.25
$ws.Cells.Item(91, 9).Value = 2400
$ws.Cells.Item(91, 10).Value = 3228.3333
$ws.Cells.Item(91, 11).Value = 2400
$ws.Cells.Item(91, 12).Value = 3228.3333
$ws.Cells.Item(91, 13).Value = -996
$ws.Cells.Item(91, 14).Value = -6036.3333

$ws.Cells.Item(113, 8).Value = 40000
$ws.Cells.Item(113, 10).Value = 40000
$ws.Cells.Item(113, 12).Value = 40000
$ws.Cells.Item(113, 14).Value = -48678

$ws.Cells.Item(136, 8).Value = 3039.8262
$ws.Cells.Item(136, 9).Value = 2651
$ws.Cells.Item(136, 10).Value = 3928.5715
$ws.Cells.Item(136, 11).Value = 7953
$ws.Cells.Item(136, 12).Value = 11785.7145
$ws.Cells.Item(136, 13).Value = -5403
$ws.Cells.Item(136, 14).Value = -16885.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 58112.5
$ws.Cells.Item(86, 9).Value = 2801.375
$ws.Cells.Item(86, 10).Value = 102361.4
$ws.Cells.Item(86, 11).Value = 2801.375
$ws.Cells.Item(86, 12).Value = 102361.4
$ws.Cells.Item(86, 13).Value = -1678.375
$ws.Cells.Item(86, 14).Value = -104607.4

$ws.Cells.Item(89, 8).Value = 58112.5
$ws.Cells.Item(89, 9).Value = 2801.375
$ws.Cells.Item(89, 10).Value = 102361.4
$ws.Cells.Item(89, 11).Value = 14006.875
$ws.Cells.Item(89, 12).Value = 511807
$ws.Cells.Item(89, 13).Value = -8390.875
$ws.Cells.Item(89, 14).Value = -523039

$ws.Cells.Item(134, 8).Value = 2578.2632
$ws.Cells.Item(134, 9).Value = 2306.6924
$ws.Cells.Item(134, 10).Value = 3166.6667
$ws.Cells.Item(134, 11).Value = 6920.0772
$ws.Cells.Item(134, 12).Value = 9500.000100000001
$ws.Cells.Item(134, 13).Value = -4385.0772
$ws.Cells.Item(134, 14).Value = -14570.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3525.0715
$ws.Cells.Item(31, 9).Value = 1847.7646
$ws.Cells.Item(31, 10).Value = 6117.273
$ws.Cells.Item(31, 11).Value = 1847.7646
$ws.Cells.Item(31, 12).Value = 6117.273
$ws.Cells.Item(31, 13).Value = -1552.7646
$ws.Cells.Item(31, 14).Value = -6707.273

$ws.Cells.Item(34, 8).Value = 3525.0715
$ws.Cells.Item(34, 9).Value = 1847.7646
$ws.Cells.Item(34, 10).Value = 6117.273
$ws.Cells.Item(34, 11).Value = 1847.7646
$ws.Cells.Item(34, 12).Value = 6117.273
$ws.Cells.Item(34, 13).Value = -1645.7646
$ws.Cells.Item(34, 14).Value = -6521.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1980.3334
$ws.Cells.Item(132, 9).Value = 777
$ws.Cells.Item(132, 10).Value = 2051.1177
$ws.Cells.Item(132, 11).Value = 6993
$ws.Cells.Item(132, 12).Value = 18460.0593
$ws.Cells.Item(132, 13).Value = -4463
$ws.Cells.Item(132, 14).Value = -23520.0593

$ws.Cells.Item(133, 8).Value = 5532.231
$ws.Cells.Item(133, 9).Value = 1745.25
$ws.Cells.Item(133, 11).Value = 5235.75
$ws.Cells.Item(133, 13).Value = -175.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 42502
$ws.Cells.Item(6, 9).Value = 22008
$ws.Cells.Item(6, 10).Value = 49333.332
$ws.Cells.Item(6, 11).Value = 22008
$ws.Cells.Item(6, 12).Value = 49333.332
$ws.Cells.Item(6, 13).Value = -21895
$ws.Cells.Item(6, 14).Value = -49559.332

$ws.Cells.Item(16, 8).Value = 42502
$ws.Cells.Item(16, 9).Value = 22008
$ws.Cells.Item(16, 10).Value = 49333.332
$ws.Cells.Item(16, 11).Value = 22008
$ws.Cells.Item(16, 12).Value = 49333.332
$ws.Cells.Item(16, 13).Value = -21758
$ws.Cells.Item(16, 14).Value = -49833.332

$ws.Cells.Item(93, 8).Value = 40000
$ws.Cells.Item(93, 10).Value = 40000
$ws.Cells.Item(93, 12).Value = 40000
$ws.Cells.Item(93, 14).Value = -43744

$ws.Cells.Item(109, 8).Value = 18856.428
$ws.Cells.Item(109, 10).Value = 18856.428
$ws.Cells.Item(109, 12).Value = 18856.428
$ws.Cells.Item(109, 14).Value = -20936.428

$ws.Cells.Item(123, 8).Value = 28723.8
$ws.Cells.Item(123, 10).Value = 28723.8
$ws.Cells.Item(123, 12).Value = 28723.8
$ws.Cells.Item(123, 14).Value = -33623.8

$ws.Cells.Item(132, 8).Value = 3635.7896
$ws.Cells.Item(132, 9).Value = 2466.2
$ws.Cells.Item(132, 10).Value = 4053.5
$ws.Cells.Item(132, 11).Value = 7398.599999999999
$ws.Cells.Item(132, 12).Value = 12160.5
$ws.Cells.Item(132, 13).Value = -4868.599999999999
$ws.Cells.Item(132, 14).Value = -17220.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5018.857
$ws.Cells.Item(7, 9).Value = 4108.8
$ws.Cells.Item(7, 11).Value = 4108.8
$ws.Cells.Item(7, 13).Value = -3996.8

$ws.Cells.Item(82, 8).Value = 1914.6296
$ws.Cells.Item(82, 9).Value = 1419
$ws.Cells.Item(82, 10).Value = 2255.375
$ws.Cells.Item(82, 11).Value = 1419
$ws.Cells.Item(82, 12).Value = 2255.375
$ws.Cells.Item(82, 13).Value = -1058
$ws.Cells.Item(82, 14).Value = -2977.375

$ws.Cells.Item(85, 8).Value = 1914.6296
$ws.Cells.Item(85, 9).Value = 1419
$ws.Cells.Item(85, 10).Value = 2255.375
$ws.Cells.Item(85, 11).Value = 1419
$ws.Cells.Item(85, 12).Value = 2255.375
$ws.Cells.Item(85, 13).Value = -171
$ws.Cells.Item(85, 14).Value = -4751.375

$ws.Cells.Item(122, 8).Value = 75002160
$ws.Cells.Item(122, 9).Value = 62502240
$ws.Cells.Item(122, 10).Value = 100002000
$ws.Cells.Item(122, 11).Value = 187506720
$ws.Cells.Item(122, 12).Value = 300006000
$ws.Cells.Item(122, 13).Value = -187504270
$ws.Cells.Item(122, 14).Value = -300010900

$ws.Cells.Item(126, 8).Value = 5018.857
$ws.Cells.Item(126, 9).Value = 4108.8
$ws.Cells.Item(126, 11).Value = 12326.4
$ws.Cells.Item(126, 13).Value = -9856.400000000001

$ws.Cells.Item(132, 8).Value = 6602.385
$ws.Cells.Item(132, 9).Value = 8822.167
$ws.Cells.Item(132, 10).Value = 4699.7144
$ws.Cells.Item(132, 11).Value = 26466.501
$ws.Cells.Item(132, 12).Value = 14099.1432
$ws.Cells.Item(132, 13).Value = -23936.501
$ws.Cells.Item(132, 14).Value = -19159.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4319
$ws.Cells.Item(62, 9).Value = 3880
$ws.Cells.Item(62, 10).Value = 4758
$ws.Cells.Item(62, 11).Value = 3880
$ws.Cells.Item(62, 12).Value = 4758
$ws.Cells.Item(62, 13).Value = -3256
$ws.Cells.Item(62, 14).Value = -6006

$ws.Cells.Item(65, 8).Value = 4319
$ws.Cells.Item(65, 9).Value = 3880
$ws.Cells.Item(65, 10).Value = 4758
$ws.Cells.Item(65, 11).Value = 19400
$ws.Cells.Item(65, 12).Value = 23790
$ws.Cells.Item(65, 13).Value = -16280
$ws.Cells.Item(65, 14).Value = -30030

$ws.Cells.Item(81, 8).Value = 47330.23
$ws.Cells.Item(81, 9).Value = 76213.6
$ws.Cells.Item(81, 11).Value = 152427.2
$ws.Cells.Item(81, 13).Value = -151366.2

$ws.Cells.Item(84, 8).Value = 47330.23
$ws.Cells.Item(84, 9).Value = 76213.6
$ws.Cells.Item(84, 11).Value = 762136
$ws.Cells.Item(84, 13).Value = -756832

$ws.Cells.Item(114, 8).Value = 80000
$ws.Cells.Item(114, 10).Value = 80000
$ws.Cells.Item(114, 12).Value = 80000
$ws.Cells.Item(114, 14).Value = -88678

$ws.Cells.Item(123, 8).Value = 40650.914
$ws.Cells.Item(123, 10).Value = 40650.914
$ws.Cells.Item(123, 12).Value = 40650.914
$ws.Cells.Item(123, 14).Value = -50450.914

$ws.Cells.Item(126, 8).Value = 4893.533
$ws.Cells.Item(126, 9).Value = 4864.273
$ws.Cells.Item(126, 10).Value = 4974
$ws.Cells.Item(126, 11).Value = 14592.819
$ws.Cells.Item(126, 12).Value = 14922
$ws.Cells.Item(126, 13).Value = -12122.819
$ws.Cells.Item(126, 14).Value = -19862

$ws.Cells.Item(132, 8).Value = 3647.6897
$ws.Cells.Item(132, 9).Value = 3313.2144
$ws.Cells.Item(132, 10).Value = 3959.8667
$ws.Cells.Item(132, 11).Value = 9939.643199999999
$ws.Cells.Item(132, 12).Value = 11879.6001
$ws.Cells.Item(132, 13).Value = -7409.643199999999
$ws.Cells.Item(132, 14).Value = -16939.6001
